# Risiken.pptx edit: fill in risk-matrix cell labels (requirement / component codes)
# on the risk matrix table (slide 1, graphicFrame id=5).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tbl = $sh.Table

# Simple empty -> text cells
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "I1"
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Text = "C2"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "C3"
$tbl.Cell(3, 6).Shape.TextFrame.TextRange.Text = "FV1"
$tbl.Cell(4, 4).Shape.TextFrame.TextRange.Text = "FF2, B1"
$tbl.Cell(4, 5).Shape.TextFrame.TextRange.Text = "I2, K1"
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "FF1, B2"

# "Sehr Hohes / IT-Risiko" cell: append the new requirement codes to the
# existing second line ("IT-Risiko" -> "IT-Risiko, I3, C1, K2")
$cell = $tbl.Cell(2, 5)
$tr = $cell.Shape.TextFrame.TextRange
$sub = $tr.Characters(13, 9)
$sub.Text = "IT-Risiko, I3, C1, K2"
